$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D to Text format first so numeric-looking price strings
# (e.g. "303.40") are preserved exactly as text instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.252.85'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.604.12'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '303.40'
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").Value = '0.3768'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '52.00'
$ws.Range("E8").Value = '  +4.06%  '
$ws.Range("D9").Value = '0.3641'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").Value = '1.281'
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.08139'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '22.88'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '6.603'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '7.435'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").Value = '0.00001249'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '1.600.85'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '94.07'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").Value = '0.06946'
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '18.21'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '6.539'
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("D24").Value = '23.246.76'
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").Value = '3.057'
$ws.Range("E25").Value = '  +9.14%  '
$ws.Range("D26").Value = '2.383'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '149.96'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '5.262'
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").Value = '134.77'
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("D31").Value = '2.394'
$ws.Range("E31").Value = '  +3.53%  '
$ws.Range("D32").Value = '6.755'
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").Value = '1.781.23'
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").Value = '0.9677'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = '0.07500'
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("D36").Value = '0.02766'
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("D37").Value = '10.32'
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = '0.2544'
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '6.129'
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.08812'
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  +2.26%  '
$ws.Range("D42").Value = '0.7130'
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").Value = '12.51'
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = '15.70'
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("D45").Value = '0.6568'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '2.323'
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("D47").Value = '0.9993'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '4.015'
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("D49").Value = '132.62'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").Value = '1.209'
$ws.Range("E51").Value = '  -0.67%  '

# Restore the default (Normal) style on column D so no stray
# number-format styling is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
